$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.308.88'
$ws.Range("E2").Value = '  +0.19%  '
$ws.Range("D3").Value = '2.278.76'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.611'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.77%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -1.47%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.58'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.23%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0898'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.49%  '
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.964'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("D16").Value = '2.624.63'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '2.288.66'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").Value = '42.329.79'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.22'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.98%  '
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '257.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.12%  '
$ws.Range("E25").Value = '  -2.81%  '
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.20%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +12.71%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.89%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0845'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.128'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.48%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -3.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.47'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.83%  '
$ws.Range("E38").Value = '  -2.40%  '
$ws.Range("E39").Value = '  -1.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.98%  '
$ws.Range("E41").Value = '  +2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.65%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.01'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '68.31'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.98%  '
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '11.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").Value = '1.701.49'
$ws.Range("E47").Value = '  +6.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '109.69'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '76.64'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.15%  '
$ws.Range("E50").Value = '  -3.65%  '
